# Applies weekly update to 'Bruselas (repollito)' price records (rows 26-47).
# A new weekly record is inserted at row 26 (shifting prior rows 26-46 down to 27-47),
# and the data that used to occupy row 47 is dropped off the bottom of this block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = 45044
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 22000
$ws.Range("M26").Value = 21000
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 1400
$ws.Range("D27").Value = 44714
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 19000
$ws.Range("P27").Value = 1267
$ws.Range("D28").Value = 44754
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = 22000
$ws.Range("L28").Value = 22000
$ws.Range("M28").Value = 22000
$ws.Range("P28").Value = 1467
$ws.Range("D29").Value = 44726
$ws.Range("J29").Value = 28
$ws.Range("L29").Value = 24000
$ws.Range("M29").Value = 24000
$ws.Range("P29").Value = 1600
$ws.Range("D30").Value = 44406
$ws.Range("J30").Value = 25
$ws.Range("K30").Value = 24000
$ws.Range("L30").Value = 25000
$ws.Range("M30").Value = 24520
$ws.Range("P30").Value = 1635
$ws.Range("D31").Value = 44817
$ws.Range("J31").Value = 18
$ws.Range("K31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = 20000
$ws.Range("P31").Value = 1333
$ws.Range("D32").Value = 44432
$ws.Range("J32").Value = 34
$ws.Range("K32").Value = 24000
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 24500
$ws.Range("P32").Value = 1633
$ws.Range("D33").Value = 44707
$ws.Range("J33").Value = 30
$ws.Range("K33").Value = 26000
$ws.Range("L33").Value = 26000
$ws.Range("M33").Value = 26000
$ws.Range("P33").Value = 1733
$ws.Range("D34").Value = 44790
$ws.Range("J34").Value = 36
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = 20000
$ws.Range("P34").Value = 1333
$ws.Range("D35").Value = 44460
$ws.Range("J35").Value = 25
$ws.Range("K35").Value = 24000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = 24480
$ws.Range("P35").Value = 1632
$ws.Range("D36").Value = 44397
$ws.Range("K36").Value = 23000
$ws.Range("L36").Value = 24000
$ws.Range("M36").Value = 23500
$ws.Range("P36").Value = 1567
$ws.Range("D37").Value = 44769
$ws.Range("K37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = 20000
$ws.Range("P37").Value = 1333
$ws.Range("D38").Value = 44411
$ws.Range("J38").Value = 34
$ws.Range("K38").Value = 25000
$ws.Range("L38").Value = 26000
$ws.Range("M38").Value = 25500
$ws.Range("P38").Value = 1700
$ws.Range("D39").Value = 44455
$ws.Range("J39").Value = 18
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = 24500
$ws.Range("P39").Value = 1633
$ws.Range("D40").Value = 44708
$ws.Range("J40").Value = 25
$ws.Range("K40").Value = 26000
$ws.Range("L40").Value = 26000
$ws.Range("M40").Value = 26000
$ws.Range("P40").Value = 1733
$ws.Range("D41").Value = 44811
$ws.Range("J41").Value = 18
$ws.Range("K41").Value = 20000
$ws.Range("M41").Value = 20000
$ws.Range("P41").Value = 1333
$ws.Range("D42").Value = 44831
$ws.Range("J42").Value = 20
$ws.Range("K42").Value = 19000
$ws.Range("M42").Value = 19500
$ws.Range("P42").Value = 1300
$ws.Range("D43").Value = 44763
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = 20000
$ws.Range("P43").Value = 1333
$ws.Range("D44").Value = 44677
$ws.Range("K44").Value = 25000
$ws.Range("L44").Value = 26000
$ws.Range("M44").Value = 25500
$ws.Range("P44").Value = 1700
$ws.Range("D45").Value = 44749
$ws.Range("J45").Value = 34
$ws.Range("K45").Value = 18000
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = 18000
$ws.Range("P45").Value = 1200
$ws.Range("D46").Value = 44784
$ws.Range("J46").Value = 28
$ws.Range("K46").Value = 20000
$ws.Range("L46").Value = 21000
$ws.Range("M46").Value = 20357
$ws.Range("P46").Value = 1357
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 25
$ws.Range("K47").Value = 14000
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = 14480
$ws.Range("P47").Value = 965
